# Waardelijst voor Reflevel toegevoegd
# Adds a new "RefLevel" worksheet (a BRO GLD value list) as the last sheet
# in the workbook, with the standard 4-column Codes/IsImbro/IsImbroA/Description
# header layout used by the other value-list sheets, and a single "NAP" entry.

$wb = $excel.ActiveWorkbook

# Add the new sheet after the current last sheet so it lands at the end of
# the tab strip (Worksheets.Add() defaults to inserting before the active
# sheet, which is not what we want here).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "RefLevel"

# Header rows (match the layout of the other BRO GLD value-list sheets).
$ws.Range("A1").Value = "Codes"
$ws.Range("B1").Value = "IsImbro"
$ws.Range("C1").Value = "IsImbroA"
$ws.Range("D1").Value = "Description"

$ws.Range("A2").Value = "[String]"
$ws.Range("B2").Value = "[Boolean]"
$ws.Range("C2").Value = "[Boolean]"
$ws.Range("D2").Value = "[String]"

# Single value-list entry for RefLevel.
$ws.Range("A3").Value = "NAP"
$ws.Range("B3").Value = $false
$ws.Range("C3").Value = $false
$ws.Range("D3").Value = "Referentie voor hoogtemetingen"

# Make the new sheet the active / selected tab, with D4 as the active cell
# (just past the last used row), mirroring the saved state of the sheet.
$ws.Activate()
$ws.Range("D4").Select()
